$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.064.44'
$ws.Range("E2").Value = '  -1.85%  '
$ws.Range("D3").Value = '2.294.05'
$ws.Range("E3").Value = '  -2.63%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'313.71"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.31%  '
$ws.Range("D6").Value = "'104.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.87%  '
$ws.Range("D7").Value = "'0.626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.09%  '
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").Value = "'0.608"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").Value = "'39.97"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.77%  '
$ws.Range("D11").Value = "'0.0910"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.78%  '
$ws.Range("D12").Value = "'8.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").Value = "'0.107"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("D14").Value = "'0.973"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.34%  '
$ws.Range("D15").Value = "'15.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -4.45%  '
$ws.Range("D16").Value = '2.638.46'
$ws.Range("E16").Value = '  -2.52%  '
$ws.Range("D17").Value = '2.267.86'
$ws.Range("E17").Value = '  -4.18%  '
$ws.Range("D18").Value = '41.966.77'
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("D19").Value = "'7.60"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.86%  '
$ws.Range("D21").Value = "'72.75"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -5.91%  '
$ws.Range("D22").Value = "'3.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.21%  '
$ws.Range("D23").Value = "'258.37"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("D24").Value = "'2.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.92%  '
$ws.Range("D25").Value = "'9.87"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.62%  '
$ws.Range("D26").Value = "'1.01"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.60%  '
$ws.Range("D27").Value = "'10.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.39%  '
$ws.Range("E28").Value = '  +2.60%  '
$ws.Range("D29").Value = "'22.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.71%  '
$ws.Range("D30").Value = "'35.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.84%  '
$ws.Range("D31").Value = "'163.81"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -6.36%  '
$ws.Range("D32").Value = "'0.0884"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.04%  '
$ws.Range("E33").Value = '  -2.72%  '
$ws.Range("D34").Value = "'5.86"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("E35").Value = '  -0.76%  '
$ws.Range("E36").Value = '  +4.77%  '
$ws.Range("D37").Value = "'4.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.53%  '
$ws.Range("E38").Value = '  +8.20%  '
$ws.Range("E39").Value = '  -2.31%  '
$ws.Range("D40").Value = "'3.62"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.89%  '
$ws.Range("D41").Value = "'100.06"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +17.25%  '
$ws.Range("E42").Value = '  +0.53%  '
$ws.Range("D43").Value = "'70.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.55%  '
$ws.Range("D44").Value = "'0.227"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.75%  '
$ws.Range("E45").Value = '  +0.21%  '
$ws.Range("D46").Value = "'12.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.50%  '
$ws.Range("D47").Value = "'114.01"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.80%  '
$ws.Range("D48").Value = "'78.14"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +6.40%  '
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("D50").Value = "'5.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -4.05%  '
$ws.Range("D51").Value = "'1.28"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.97%  '
